# Fix DB Save Vendings
# Updates the "Machines" sheet:
#   - Row 3 (ID=2): Opera/2019/213  -> Opel/1243/14   (Country stays "France")
#   - Adds row 4 (ID=3): Opera__XAML__ / 4312 / 4123 / g53
#   - Adds row 5 (ID=4): OPLA / 1423 / 312 / Germany

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Machines")

# --- Row 3: correct the vendor info (keep it stored as text, like the rest
#     of the sheet, rather than letting Excel auto-convert the numeric-
#     looking values into numbers) ---
$ws.Range("C3:D3").NumberFormat = "@"
$ws.Range("B3").Value = "Opel"
$ws.Range("C3").Value = "1243"
$ws.Range("D3").Value = "14"

# --- Row 4 (new) ---
$ws.Range("A4").NumberFormat = "@"
$ws.Range("C4:D4").NumberFormat = "@"
$ws.Range("A4").Value = "3"
$ws.Range("B4").Value = "Opera__XAML__"
$ws.Range("C4").Value = "4312"
$ws.Range("D4").Value = "4123"
$ws.Range("E4").Value = "g53"

# --- Row 5 (new) ---
$ws.Range("A5").NumberFormat = "@"
$ws.Range("C5:D5").NumberFormat = "@"
$ws.Range("A5").Value = "4"
$ws.Range("B5").Value = "OPLA"
$ws.Range("C5").Value = "1423"
$ws.Range("D5").Value = "312"
$ws.Range("E5").Value = "Germany"
